$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lineup_Template")

# Normalize the redundant "fill" formatting on A14:A17 so it matches the
# (border + vertical-center, no fill) styling used by the other label cells.
$ws.Range("A14:A17").Interior.Pattern = -4142

# Fill out the lineup form (values chosen from the HIDDEN sheet's validation lists)
$ws.Range("B3").Value = 1
$ws.Range("B5").Value = "Sandy T"
$ws.Range("B6").Value = "Mike K"
$ws.Range("B7").Value = "Mike K"
$ws.Range("B8").Value = "John J"
$ws.Range("B9").Value = "Adam A"
$ws.Range("B10").Value = "Brooks K"
$ws.Range("B11").Value = "Brooks K"
$ws.Range("B12").Value = "Scott S"
$ws.Range("B13").Value = "Tiger W"
$ws.Range("B14").Value = "Albert O"
$ws.Range("B15").Value = "Christy J"
$ws.Range("B16").Value = "Albert O"
$ws.Range("B17").Value = "Albert O"
$ws.Range("B18").Value = "Houston"

# Move the active selection to the last cell that was edited
$ws.Range("B18").Select()
